# fall 22 week 14 day-after
# Adds a new "Week 48" column (AW) to the Inning Counts sheet and fills
# in the handful of player values recorded for that week.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header (shared string "Week 48", appended after "Week 47")
$ws.Range("AW1").Value = "Week 48"

# New week's data for the players that have a recorded value
$ws.Range("AW2").Value = 3.5    # Scott Foxley-Berry
$ws.Range("AW5").Value = 10     # Jason Bohrer
$ws.Range("AW6").Value = 5      # Dan Aquino
$ws.Range("AW7").Value = 5.5    # Jason Liess
$ws.Range("AW8").Value = 10     # Kim Quan

# Keep the window/selection state close to what was recorded after the edit
[void]$ws.Range("AV12").Select()
